$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 284 (shifts existing rows 284:360 down to 285:361,
# carrying their values/formatting with them).
$ws.Range("A284:R284").EntireRow.Insert()

# Populate the newly inserted row 284 with the new data record.
$ws.Range("A284").Value = 10
$ws.Range("B284").Value = "Vega Modelo de Temuco"
$ws.Range("C284").Value = "La Araucanía"
$ws.Range("D284").Value = 45204
$ws.Range("E284").Value = 9
$ws.Range("F284").Value = 100112013
$ws.Range("G284").Value = "Alcachofa"
$ws.Range("H284").Value = "Española"
$ws.Range("I284").Value = "Primera"
$ws.Range("J284").Value = 185
$ws.Range("K284").Value = 12000
$ws.Range("L284").Value = 12000
$ws.Range("M284").Value = 12000
$ws.Range("N284").Value = "$/caja 30 unidades"
$ws.Range("O284").Value = "Provincia de Limarí"
$ws.Range("P284").Value = 400
$ws.Range("Q284").Value = 30
$ws.Range("R284").Value = "Hortaliza"
